$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("projected births")

# Remove the bad formula reference to "current year" (2016) births in row 2
# (it pulled from demographics!$B$3). Deleting the entire row shifts all the
# subsequent hardcoded year/births rows up by one, so the sheet now starts
# at 2017 and ends at 2030 with plain numeric values throughout.
$ws.Rows.Item(2).Delete()
